$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 216, pushing the existing rows 216-328
# down to 217-329 (dimension grows from A1:R328 to A1:R329).
$ws.Rows(216).Insert()

# Populate the newly inserted row 216 with the new daily record.
$ws.Cells.Item(216, 1).Value  = 4
$ws.Cells.Item(216, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(216, 3).Value  = "Los Lagos"
$ws.Cells.Item(216, 4).Value  = 44806
$ws.Cells.Item(216, 5).Value  = 10
$ws.Cells.Item(216, 6).Value  = 100112037
$ws.Cells.Item(216, 7).Value  = "Cebollín"
$ws.Cells.Item(216, 8).Value  = "Sin especificar"
$ws.Cells.Item(216, 9).Value  = "Primera"
$ws.Cells.Item(216, 10).Value = 180
$ws.Cells.Item(216, 11).Value = 9000
$ws.Cells.Item(216, 12).Value = 10000
$ws.Cells.Item(216, 13).Value = 9500
$ws.Cells.Item(216, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(216, 15).Value = "Región Metropolitana"
$ws.Cells.Item(216, 16).Value = 264
$ws.Cells.Item(216, 17).Value = 36
$ws.Cells.Item(216, 18).Value = "Hortaliza"
